# EA PB Awards tables - rebuild update
# - Rename the "Bucket" category labels in column A from singular to plural form:
#     Sprint -> Sprints, Jump -> Jumps, Throw -> Throws
# - Reset the sheet view (scroll position / active selection) back to the top
#   of the sheet instead of where the previous editor had left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Constants for Excel's Find/Replace (xlWhole = 1, xlByRows = 1, xlPart = 2)
$xlWhole = 1

# Rename category bucket labels (column A) - exact whole-cell matches only,
# so this can't accidentally touch any other text in the sheet.
$ws.Cells.Replace("Sprint", "Sprints", $xlWhole) | Out-Null
$ws.Cells.Replace("Throw", "Throws", $xlWhole) | Out-Null
$ws.Cells.Replace("Jump", "Jumps", $xlWhole) | Out-Null

# Restore the sheet view to the top-left and move the active selection to O22.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("O22").Select()

$wb.Save()
